$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.237.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.40%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.851.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.31%  "
# Row 4
$ws.Range("E4").Value = "  -0.07%  "
# Row 5
$ws.Range("B5").Value = "XRP"
$ws.Range("C5").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.6957"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.00%  "
# Row 6
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "238.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.67%  "
# Row 7
$ws.Range("E7").Value = "  -0.14%  "
# Row 8
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3076"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.23%  "
# Row 9
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07612"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.51%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.18%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08090"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.80%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7234"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.86%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.836.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.55%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.202"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.99%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.63%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.047.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.07%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.879"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.74%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "242.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.37%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007743"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.38%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.50%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.05%  "
# Row 22
$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "
# Row 23
$ws.Range("B23").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C23").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.052.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.25%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.606"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.18%  "
# Row 25
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "163.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.19%  "
# Row 26
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.069"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.12%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1460"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.87%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.22%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.936"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.91%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.400"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.19%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.504"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.08%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.434"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.87%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.038"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.36%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05278"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.02%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.192"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.48%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7130"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.65%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.002"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.22%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.664"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.54%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01860"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.79%  "
# Row 40
$ws.Range("E40").Value = "  -2.12%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9355"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.32%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4302"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.58%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.882"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.09%  "
# Row 44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.044.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.15%  "
# Row 45
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "69.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.88%  "
# Row 46
$ws.Range("E46").Value = "  -0.19%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.47%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.263"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.55%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.740"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.18%  "
# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.275"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.25%  "
# Row 51
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.972.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.26%  "
